# Update NATMI TPM edge-weight table with new TPM values.
# - Adds a new target cluster "Resolving-Mac" (2 new rows: FAPs->Resolving-Mac, MuSCs->Resolving-Mac)
# - Recomputes receptor/edge-weight derived columns (K:T) for every Sending/Target cluster pair
# - Re-orders rows so they are grouped by Sending cluster (FAPs, then MuSCs) and, within each
#   group, by Target cluster (ECs, FAPs, MuSCs, Resolving-Mac)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2: FAPs -> ECs (receptor/edge columns recomputed) ----
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.3643233333333333
$ws.Range("N2").Value = 1.09297
$ws.Range("O2").Value = 0.09146135066097912
$ws.Range("P2").Value = 0.09146135066097912
$ws.Range("Q2").Value = 3.801438312011111
$ws.Range("R2").Value = 34.2129448081
$ws.Range("S2").Value = 0.08881895887732785
$ws.Range("T2").Value = 0.08881895887732787

# ---- Row 3: FAPs -> FAPs (specificity columns recomputed) ----
$ws.Range("O3").Value = 0.2211540931751474
$ws.Range("P3").Value = 0.2211540931751474
$ws.Range("S3").Value = 0.2147647740310101
$ws.Range("T3").Value = 0.2147647740310101

# ---- Row 4: FAPs -> MuSCs (receptor/edge columns recomputed) ----
$ws.Range("M4").Value = 2.721212
$ws.Range("N4").Value = 8.163636
$ws.Range("O4").Value = 0.6831451685449673
$ws.Range("P4").Value = 0.6831451685449673
$ws.Range("Q4").Value = 28.39378816958666
$ws.Range("R4").Value = 255.54409352628
$ws.Range("S4").Value = 0.6634085566607256
$ws.Range("T4").Value = 0.6634085566607257

# ---- Row 5: now FAPs -> Resolving-Mac (new pairing takes this row's place) ----
$ws.Range("A5").Value = "FAPs"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 10.43424333333333
$ws.Range("H5").Value = 31.30273
$ws.Range("I5").Value = 0.9711091978791583
$ws.Range("J5").Value = 0.9711091978791584
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.016887
$ws.Range("N5").Value = 0.050661
$ws.Range("O5").Value = 0.004239387618906157
$ws.Range("P5").Value = 0.004239387618906157
$ws.Range("Q5").Value = 0.17620306717
$ws.Range("R5").Value = 1.58582760453
$ws.Range("S5").Value = 0.004116908310094793
$ws.Range("T5").Value = 0.004116908310094793

# ---- Row 6: now MuSCs -> ECs ----
$ws.Range("D6").Value = "ECs"
$ws.Range("M6").Value = 0.3643233333333333
$ws.Range("N6").Value = 1.09297
$ws.Range("O6").Value = 0.09146135066097912
$ws.Range("P6").Value = 0.09146135066097912
$ws.Range("Q6").Value = 0.11309397778
$ws.Range("R6").Value = 1.01784580002
$ws.Range("S6").Value = 0.002642391783651253
$ws.Range("T6").Value = 0.002642391783651254

# ---- Row 7: now MuSCs -> FAPs ----
$ws.Range("D7").Value = "FAPs"
$ws.Range("M7").Value = 0.8809360000000001
$ws.Range("N7").Value = 2.642808
$ws.Range("O7").Value = 0.2211540931751474
$ws.Range("P7").Value = 0.2211540931751474
$ws.Range("Q7").Value = 0.273461914992
$ws.Range("R7").Value = 2.461157234928
$ws.Range("S7").Value = 0.006389319144137352
$ws.Range("T7").Value = 0.006389319144137353

# ---- Row 8 (new): MuSCs -> MuSCs ----
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Wnt5a"
$ws.Range("C8").Value = "Fzd3"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.310422
$ws.Range("H8").Value = 0.9312659999999999
$ws.Range("I8").Value = 0.02889080212084161
$ws.Range("J8").Value = 0.02889080212084161
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 2.721212
$ws.Range("N8").Value = 8.163636
$ws.Range("O8").Value = 0.6831451685449673
$ws.Range("P8").Value = 0.6831451685449673
$ws.Range("Q8").Value = 0.844724071464
$ws.Range("R8").Value = 7.602516643176
$ws.Range("S8").Value = 0.01973661188424164
$ws.Range("T8").Value = 0.01973661188424164

# ---- Row 9 (new): MuSCs -> Resolving-Mac ----
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Wnt5a"
$ws.Range("C9").Value = "Fzd3"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.310422
$ws.Range("H9").Value = 0.9312659999999999
$ws.Range("I9").Value = 0.02889080212084161
$ws.Range("J9").Value = 0.02889080212084161
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.016887
$ws.Range("N9").Value = 0.050661
$ws.Range("O9").Value = 0.004239387618906157
$ws.Range("P9").Value = 0.004239387618906157
$ws.Range("Q9").Value = 0.005242096314
$ws.Range("R9").Value = 0.047178866826
$ws.Range("S9").Value = 0.0001224793088113637
$ws.Range("T9").Value = 0.0001224793088113637
